$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename A16/A17 back from onEnter/onLeave to enterCB/leaveCB
$ws.Range("A16").Value = "enterCB"
$ws.Range("A17").Value = "leaveCB"

# Delete row 18 (the "onChange" row added by the reverted commit) and shift rows up
$ws.Rows.Item(18).Delete()

# Update selection to match target state
$ws.Range("L15").Select()
